$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Execute" flag (column B) from "No" to "Yes" for the rows that
# were turned on in this regression-testing pass.
$executeRows = @(12, 22, 23, 24, 25, 26, 27, 28, 32, 34, 35, 36)
foreach ($r in $executeRows) {
    $ws.Cells.Item($r, 2).Value = "Yes"
}

# Fix the parameter name typo in the testAddSignetAccountInvalidAmount test
# data (-amount -> -pamount).
$ws.Range("I19").Value = "coyni.admin.tests.CoyniPortalTest,
testAddSignetAccountInvalidAmount,
-pamount,
-pdescription,
-perrMessage
"

# Reflect the cursor/selection position left behind by the author when the
# workbook was saved.
[void]$ws.Range("B36").Select()
